# Generate Report for Handback
#
# 1. Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn / de-de status columns) and on the
#    per-language "Status" column of the zh-cn / de-de sheets.
# 2. Each per-language sheet gains two new columns: F = "Latest Target
#    File" and G = "Latest Handback File" - mirroring the source file
#    (column A) and the handoff target file (column D), each as a
#    hyperlinked file name.
# 3. The "Latest Handback DateTime" column (H) is stamped with the
#    handback timestamp (different per language sheet).

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both language-status columns (B = zh-cn, C = de-de)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack

# Latest Target File / Latest Handback File - row 2 (106d6da0... file)
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", "106d6da0-5c15-4669-815c-ad923b15a0fc.md")
$zhcn.Range("F2").Style = "HyperLink"

$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2260a60e2799454237861fb46b4fd2470a45ff2a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.zh-cn.xlf", "", "", "106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.zh-cn.xlf")
$zhcn.Range("G2").Style = "HyperLink"

# Latest Target File / Latest Handback File - row 3 (9e74227a... file)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", "9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md")
$zhcn.Range("F3").Style = "HyperLink"

$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2260a60e2799454237861fb46b4fd2470a45ff2a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.zh-cn.xlf", "", "", "9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.zh-cn.xlf")
$zhcn.Range("G3").Style = "HyperLink"

# Latest Handback DateTime
$zhcn.Range("H2").Value = "2016-03-20 10:48:53"
$zhcn.Range("H3").Value = "2016-03-20 10:48:53"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack

# Latest Target File / Latest Handback File - row 2 (106d6da0... file)
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/106d6da0-5c15-4669-815c-ad923b15a0fc.md", "", "", "106d6da0-5c15-4669-815c-ad923b15a0fc.md")
$dede.Range("F2").Style = "HyperLink"

$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5e4094ed9d048769903debf9fc1f9097c5a43b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.de-de.xlf", "", "", "106d6da0-5c15-4669-815c-ad923b15a0fc.2298611cd95f265d4cb02f723b66d1fd51448994.de-de.xlf")
$dede.Range("G2").Style = "HyperLink"

# Latest Target File / Latest Handback File - row 3 (9e74227a... file)
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed255ab16f76fe2f1191fb19786fa4f2f5a8af75/e2e/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md", "", "", "9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.md")
$dede.Range("F3").Style = "HyperLink"

$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e5e4094ed9d048769903debf9fc1f9097c5a43b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.de-de.xlf", "", "", "9e74227a-ff8b-4bb4-ba8b-887cae09bdf2.bb764def15a45df9d5253695e7594fdd7ff7de01.de-de.xlf")
$dede.Range("G3").Style = "HyperLink"

# Latest Handback DateTime
$dede.Range("H2").Value = "2016-03-20 10:48:58"
$dede.Range("H3").Value = "2016-03-20 10:48:58"
